$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.998.75'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '1.642.08'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '213.65'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '23.66'
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").Value = '0.0878'
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").Value = '1.875.38'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.641.82'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("E14").Value = '  +1.33%  '
$ws.Range("D15").Value = '0.576'
$ws.Range("E15").Value = '  +4.15%  '
$ws.Range("D16").Value = '66.01'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Value = '27.990.10'
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = '232.63'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = '7.63'
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '10.73'
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").Value = '151.65'
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  +1.34%  '
$ws.Range("D27").Value = '15.74'
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").Value = '3.34'
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").Value = '1.412.39'
$ws.Range("E34").Value = '  -4.39%  '
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("D37").Value = '0.890'
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").Value = '0.919'
$ws.Range("E40").Value = '  -4.65%  '
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +6.87%  '
$ws.Range("E44").Value = '  -2.19%  '
$ws.Range("D45").Value = '5.47'
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = '1.783.90'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").Value = '88.12'
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '7.61'
$ws.Range("E51").Value = '  -1.42%  '
